$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.015347361564636
$ws.Range("B1").Value = 1.643430590629578
$ws.Range("C1").Value = 5.102684497833252
$ws.Range("D1").Value = 1.581733584403992
$ws.Range("E1").Value = 0.3055053353309631
